$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.799.67'
$ws.Range("E2").Value = '  +6.05%  '
$ws.Range("D3").Value = '3.533.85'
$ws.Range("E3").Value = '  +9.53%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '566.56'
$ws.Range("E5").Value = '  +7.39%  '
$ws.Range("D6").Value = '188.21'
$ws.Range("E6").Value = '  +10.10%  '
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  +4.77%  '
$ws.Range("D8").Value = '3.527.17'
$ws.Range("E8").Value = '  +9.33%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = '0.634'
$ws.Range("E10").Value = '  +5.06%  '
$ws.Range("E11").Value = '  +14.38%  '
$ws.Range("D12").Value = '54.66'
$ws.Range("E12").Value = '  +3.41%  '
$ws.Range("D13").Value = '0.0000269'
$ws.Range("E13").Value = '  +6.48%  '
$ws.Range("D14").Value = '9.41'
$ws.Range("E14").Value = '  +3.42%  '
$ws.Range("D15").Value = '4.094.52'
$ws.Range("E15").Value = '  +9.38%  '
$ws.Range("D16").Value = '3.530.87'
$ws.Range("E16").Value = '  +9.22%  '
$ws.Range("E17").Value = '  +4.74%  '
$ws.Range("D18").Value = '66.801.10'
$ws.Range("E18").Value = '  +6.18%  '
$ws.Range("D19").Value = '18.23'
$ws.Range("E19").Value = '  +6.32%  '
$ws.Range("D20").Value = '12.03'
$ws.Range("E20").Value = '  +9.05%  '
$ws.Range("D21").Value = '0.997'
$ws.Range("E21").Value = '  +3.38%  '
$ws.Range("D22").Value = '427.94'
$ws.Range("E22").Value = '  +17.14%  '
$ws.Range("D23").Value = '4.17'
$ws.Range("E23").Value = '  +11.81%  '
$ws.Range("D24").Value = '85.22'
$ws.Range("E24").Value = '  +5.38%  '
$ws.Range("D25").Value = '4.13'
$ws.Range("E25").Value = '  +3.52%  '
$ws.Range("D26").Value = '11.11'
$ws.Range("E26").Value = '  -0.30%  '
$ws.Range("D27").Value = '2.90'
$ws.Range("E27").Value = '  +10.04%  '
$ws.Range("D28").Value = '12.23'
$ws.Range("E28").Value = '  +9.11%  '
$ws.Range("D29").Value = '9.23'
$ws.Range("E29").Value = '  +13.00%  '
$ws.Range("D30").Value = '30.39'
$ws.Range("E30").Value = '  +7.15%  '
$ws.Range("D31").Value = '643.00'
$ws.Range("E31").Value = '  +1.83%  '
$ws.Range("D32").Value = '6.60'
$ws.Range("E32").Value = '  +2.73%  '
$ws.Range("D33").Value = '11.74'
$ws.Range("E33").Value = '  +5.19%  '
$ws.Range("E34").Value = '  +6.35%  '
$ws.Range("D35").Value = '59.82'
$ws.Range("E35").Value = '  +6.11%  '
$ws.Range("E36").Value = '  +5.88%  '
$ws.Range("E37").Value = '  +20.39%  '
$ws.Range("D38").Value = '0.0₃0811'
$ws.Range("E38").Value = '  +14.23%  '
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("D40").Value = '0.391'
$ws.Range("E40").Value = '  +4.27%  '
$ws.Range("D41").Value = '3.35'
$ws.Range("E41").Value = '  +14.02%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").Value = '3.038.18'
$ws.Range("E43").Value = '  +5.96%  '
$ws.Range("D44").Value = '2.66'
$ws.Range("E44").Value = '  +4.94%  '
$ws.Range("D45").Value = '2.89'
$ws.Range("E45").Value = '  +11.98%  '
$ws.Range("D46").Value = '3.35'
$ws.Range("E46").Value = '  +8.25%  '
$ws.Range("E47").Value = '  +6.86%  '
$ws.Range("D48").Value = '2.75'
$ws.Range("E48").Value = '  +2.41%  '
$ws.Range("E49").Value = '  +5.95%  '
$ws.Range("D50").Value = '143.61'
$ws.Range("E50").Value = '  +7.70%  '
$ws.Range("E51").Value = '  +11.25%  '
